$d = $word.ActiveDocument

$pairs = @(
    @("15×92=1380", "61×13=793"),
    @("59×78=4602", "20×31=620"),
    @("99×39=3861", "71×80=5680"),
    @("86×45=3870", "16×60=960"),
    @("67×84=5628", "12×64=768"),
    @("56×71=3976", "70×87=6090"),
    @("47×69=3243", "29×47=1363"),
    @("55×43=2365", "87×71=6177"),
    @("90×27=2430", "23×49=1127"),
    @("39×96=3744", "33×77=2541"),
    @("82×93=7626", "15×96=1440"),
    @("68×25=1700", "34×42=1428"),
    @("71×28=1988", "95×82=7790"),
    @("32×16=512",  "41×67=2747"),
    @("98×33=3234", "18×33=594"),
    @("70×23=1610", "42×55=2310"),
    @("85×96=8160", "40×63=2520"),
    @("79×25=1975", "79×51=4029"),
    @("64×78=4992", "29×65=1885"),
    @("21×33=693",  "69×49=3381"),
    @("70×77=5390", "19×67=1273"),
    @("98×82=8036", "86×69=5934"),
    @("53×63=3339", "57×98=5586"),
    @("48×83=3984", "18×87=1566"),
    @("53×49=2597", "15×96=1440")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
